$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.617.31'
$ws.Range('E2').Value = '  +0.70%  '
$ws.Range('D3').Value = '1.565.73'
$ws.Range('E3').Value = '  -0.49%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '''210.53'
$ws.Range('E5').Value = '  -0.71%  '
$ws.Range('D6').Value = '''0.486'
$ws.Range('E6').Value = '  -1.05%  '
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Value = '''24.85'
$ws.Range('E8').Value = '  +4.76%  '
$ws.Range('E9').Value = '  -0.84%  '
$ws.Range('E10').Value = '  -0.43%  '
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('D12').Value = '1.789.64'
$ws.Range('E12').Value = '  -0.48%  '
$ws.Range('D13').Value = '1.566.67'
$ws.Range('E13').Value = '  -0.42%  '
$ws.Range('D14').Value = '28.650.34'
$ws.Range('E14').Value = '  +0.82%  '
$ws.Range('E15').Value = '  -0.52%  '
$ws.Range('E16').Value = '  -1.64%  '
$ws.Range('D17').Value = '''61.33'
$ws.Range('E17').Value = '  -0.55%  '
$ws.Range('D18').Value = '''231.74'
$ws.Range('E18').Value = '  +0.76%  '
$ws.Range('D19').Value = '''7.37'
$ws.Range('E19').Value = '  -0.33%  '
$ws.Range('D20').Value = '0.0₃0675'
$ws.Range('E20').Value = '  -1.30%  '
$ws.Range('E21').Value = '  -0.12%  '
$ws.Range('E22').Value = '  -1.23%  '
$ws.Range('D23').Value = '''8.96'
$ws.Range('E23').Value = '  -0.60%  '
$ws.Range('D24').Value = '''2.08'
$ws.Range('E24').Value = '  +1.78%  '
$ws.Range('D25').Value = '''150.38'
$ws.Range('E25').Value = '  -0.67%  '
$ws.Range('E26').Value = '  -0.92%  '
$ws.Range('E27').Value = '  -0.21%  '
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('E29').Value = '  -2.35%  '
$ws.Range('E30').Value = '  -4.80%  '
$ws.Range('D31').Value = '''1.07'
$ws.Range('E31').Value = '  -1.35%  '
$ws.Range('E32').Value = '  -0.66%  '
$ws.Range('D33').Value = '1.390.92'
$ws.Range('E33').Value = '  +0.66%  '
$ws.Range('D35').Value = '''1.03'
$ws.Range('E35').Value = '  -2.93%  '
$ws.Range('E36').Value = '  -2.20%  '
$ws.Range('E37').Value = '  +1.01%  '
$ws.Range('E38').Value = '  -2.71%  '
$ws.Range('D39').Value = '''0.0161'
$ws.Range('E39').Value = '  -1.10%  '
$ws.Range('D40').Value = '''1.94'
$ws.Range('E40').Value = '  +2.27%  '
$ws.Range('D41').Value = '''0.518'
$ws.Range('E41').Value = '  -0.46%  '
$ws.Range('D42').Value = '''1.00'
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D43').Value = '''0.772'
$ws.Range('E43').Value = '  -2.00%  '
$ws.Range('E44').Value = '  -2.78%  '
$ws.Range('D45').Value = '''63.74'
$ws.Range('E45').Value = '  +2.19%  '
$ws.Range('D46').Value = '''5.22'
$ws.Range('E46').Value = '  -2.57%  '
$ws.Range('D47').Value = '1.701.69'
$ws.Range('E47').Value = '  -0.47%  '
$ws.Range('D48').Value = '''0.870'
$ws.Range('E48').Value = '  -5.59%  '
$ws.Range('D49').Value = '''85.19'
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('D50').Value = '''43.20'
$ws.Range('E50').Value = '  +6.06%  '
$ws.Range('D51').Value = '0.0₆0102'
$ws.Range('E51').Value = '  +1.84%  '
